$d = $word.ActiveDocument

# Locate the short author byline paragraph ("Edison Achalma" styled as "Author"),
# which appears right under the "Editar: Editar" title heading. This is distinct
# from the later, longer "Nota de Autores" section that also contains this name.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Text.Trim() -eq "Edison Achalma") -and ($p.Style.NameLocal -eq "Author")) {
        $target = $p
        break
    }
}

$origIndex = $target.Index

# Insert a brand-new paragraph break right after the "Edison Achalma" paragraph.
$insertionPoint = $d.Range($target.Range.End, $target.Range.End)
$insertionPoint.InsertParagraphAfter()

# The newly created (still empty) paragraph now sits immediately after the
# original "Edison Achalma" paragraph; give it the Author style and its text.
$newPara = $d.Paragraphs.Item($origIndex + 1)
$newPara.Style = $d.Styles.Item("Author")
$newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
